$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
Write-Output $win.Left
Write-Output $win.Top
Write-Output $win.Width
Write-Output $win.Height
$win.Left = 1815
$win.Top = 1815
$win.Width = 21600
$win.Height = 10920
Write-Output "set done"
Write-Output $win.Left
